$wb = $excel.ActiveWorkbook

$newTime = "02:48:47"

# --- Sheet 1: LP1912 ---
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: $newTime"
$ws1.Range("A3").Value = "Total filas: 4"

$ws1.Range("A6").Value = $newTime
$ws1.Range("D6").Value = 10

$ws1.Range("A7").Value = $newTime
$ws1.Range("D7").Value = 60

$ws1.Range("A8").Value = $newTime
$ws1.Range("D8").Value = 73

$ws1.Range("A9").Value = $newTime
$ws1.Range("B9").Value = "04:45"
$ws1.Range("C9").Value = "215A_EL PATO"
$ws1.Range("D9").Value = 117
$ws1.Range("E9").Value = "LP1912"

# --- Sheet 2: LP1912-215 ---
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: $newTime"
$ws2.Range("A3").Value = "Total filas: 2"

$ws2.Range("A6").Value = $newTime
$ws2.Range("D6").Value = 10

$ws2.Range("A7").Value = $newTime
$ws2.Range("B7").Value = "04:45"
$ws2.Range("C7").Value = "215A_EL PATO"
$ws2.Range("D7").Value = 117
$ws2.Range("E7").Value = "LP1912"

# --- Sheet 3: 6203-6173 ---
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: $newTime"
